# Update cryptocurrency price (D) and volume-change (E) columns
# to refreshed values scraped for this run, preserving the original
# plain-text cell type (these columns are formatted/text values,
# not numeric cells) by forcing Text number format before assignment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.414.27"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.875.62"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.75"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4797"
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2831"
$ws.Range("E8").Value = "  -2.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06536"
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.876.51"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07467"
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.71"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.103"
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.36"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6634"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.392.30"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.36"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007631"
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.120.18"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.317"
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "219.05"
$ws.Range("E23").Value = "  +13.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.241"
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.360"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.32"
$ws.Range("E26").Value = "  +1.82%  "
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.983"
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.458"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09431"
$ws.Range("E30").Value = "  +3.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.316"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.049"
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05076"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.209"
$ws.Range("E34").Value = "  +5.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7535"
$ws.Range("E35").Value = "  +3.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.710"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01829"
$ws.Range("E37").Value = "  +1.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.617"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.078"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9090"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "107.00"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.910"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4297"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.461"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.73"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1285"
$ws.Range("E47").Value = "  -3.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.481"
$ws.Range("E48").Value = "  -6.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.957"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.63"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3917"
$ws.Range("E51").Value = "  +1.52%  "
